$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated average_county_temperature (I), worst_ashp_cop (N) and best_ashp_cop (O)
# values recalculated after refreshing the temperature dataset with NOAA data.

# Row 2
$ws.Range("I2").Value = 15.74228395061728
$ws.Range("N2").Value = 1.837513876759573
$ws.Range("O2").Value = 2.005936573945218

# Row 3
$ws.Range("I3").Value = 13.46442495126706
$ws.Range("N3").Value = 1.798225615362447
$ws.Range("O3").Value = 1.958604378795604

# Row 4
$ws.Range("I4").Value = 18.89814814814816
$ws.Range("N4").Value = 1.894871325212932
$ws.Range("O4").Value = 2.075424331741031

# Row 5
$ws.Range("I5").Value = 13.46442495126706
$ws.Range("N5").Value = 1.798225615362447
$ws.Range("O5").Value = 1.958604378795604

# Row 8
$ws.Range("I8").Value = 15.74228395061728
$ws.Range("N8").Value = 1.837513876759573
$ws.Range("O8").Value = 2.005936573945218

# Row 9
$ws.Range("I9").Value = 21.28240740740739
$ws.Range("N9").Value = 1.940636870984383
$ws.Range("O9").Value = 2.131200751448103

# Row 10
$ws.Range("I10").Value = 15.74228395061728
$ws.Range("N10").Value = 1.837513876759573
$ws.Range("O10").Value = 2.005936573945218
